# MVP for adding more pounds to specialty nforx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 24: mark "Has Specialty" as TRUE and fill in specialty details
$ws.Range("F24").Value = $true
$ws.Range("G24").Value = "nforx"
$ws.Range("H24").Value = 6
$ws.Range("I24").Value = 3
$ws.Range("J24").Value = 8

# Update view state to match saved selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("K24").Select()
